# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N32").Value = -5999
$ws.Range("H32").Value = 5220.8184
$ws.Range("L32").Value = 5347
$ws.Range("J32").Value = 5347
$ws.Range("M33").Value = -5555624
$ws.Range("I33").Value = 5555853
$ws.Range("K33").Value = 5555853
$ws.Range("H33").Value = 3846427.8
$ws.Range("K40").Value = 2500
$ws.Range("H40").Value = 2366.6667
$ws.Range("N40").Value = -2450
$ws.Range("M40").Value = -2325
$ws.Range("J40").Value = 2100
$ws.Range("I40").Value = 2500
$ws.Range("L40").Value = 2100
$ws.Range("H111").Value = 1913.75
$ws.Range("K111").Value = 1555.125
$ws.Range("I111").Value = 518.375
$ws.Range("M111").Value = 1511.875
$ws.Range("J113").Value = 14655.889
$ws.Range("K113").Value = 8435
$ws.Range("I113").Value = 8435
$ws.Range("N113").Value = -21163.889
$ws.Range("L113").Value = 14655.889
$ws.Range("H113").Value = 11728.412
$ws.Range("M113").Value = -5181
$ws.Range("J132").Value = 3934
$ws.Range("I132").Value = 2024.3334
$ws.Range("L132").Value = 11802
$ws.Range("H132").Value = 2569.9524
$ws.Range("M132").Value = -3543.0002
$ws.Range("K132").Value = 6073.0002
$ws.Range("N132").Value = -16862
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").Value = -5226
$ws.Range("I2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("J2").Value = 5000
$ws.Range("M2").Value = -4887
$ws.Range("K2").Value = 5000
$ws.Range("H2").Value = 5000
$ws.Range("M32").Value = -5997
$ws.Range("H32").Value = 10838
$ws.Range("I32").Value = 6284
$ws.Range("K32").Value = 6284
$ws.Range("K116").Value = 5000
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("M116").Value = -2706
$ws.Range("N116").Value = -9588
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("J130").Value = 129714.5
$ws.Range("L130").Value = 129714.5
$ws.Range("N130").Value = -139754.5
$ws.Range("H130").Value = 129714.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("H3").Value = 5000
$ws.Range("M3").Value = -4886
$ws.Range("N3").Value = -5228
$ws.Range("I3").Value = 5000
$ws.Range("H40").Value = 32222
$ws.Range("N40").Value = -32752
$ws.Range("J40").Value = 32222
$ws.Range("L40").Value = 32222
$ws.Range("M96").Value = -7536.4
$ws.Range("K96").Value = 10282.4
$ws.Range("I96").Value = 10282.4
$ws.Range("H96").Value = 14309.333
$ws.Range("K105").Value = 1294.625
$ws.Range("M105").Value = 452.375
$ws.Range("H105").Value = 1628.4117
$ws.Range("I105").Value = 1294.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 3674.75
$ws.Range("M16").Value = -3387.75
$ws.Range("H16").Value = 4284
$ws.Range("K16").Value = 3674.75
$ws.Range("H62").Value = 15183.333
$ws.Range("I62").Value = 15183.333
$ws.Range("K62").Value = 15183.333
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = -14559.333
$ws.Range("L62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("I65").Value = 15183.333
$ws.Range("H65").Value = 15183.333
$ws.Range("K65").Value = 75916.66500000001
$ws.Range("N65").ClearContents()
$ws.Range("M65").Value = -72796.66500000001
$ws.Range("K93").Value = 7146
$ws.Range("H93").Value = 7146
$ws.Range("I93").Value = 7146
$ws.Range("M93").Value = -5274
$ws.Range("H94").Value = 1770
$ws.Range("K94").Value = 1111.2
$ws.Range("M94").Value = -660.2
$ws.Range("I94").Value = 1111.2
$ws.Range("K113").Value = 3674.75
$ws.Range("I113").Value = 3674.75
$ws.Range("H113").Value = 4284
$ws.Range("M113").Value = -1504.75
$ws.Range("H134").Value = 20722.965
$ws.Range("K134").Value = 22765.092
$ws.Range("I134").Value = 7588.364
$ws.Range("M134").Value = -20230.092
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N70").Value = -17614.9995
$ws.Range("I70").Value = 4666.6665
$ws.Range("K70").Value = 13999.9995
$ws.Range("H70").Value = 5330
$ws.Range("J70").Value = 5661.6665
$ws.Range("L70").Value = 16984.9995
$ws.Range("M70").Value = -13684.9995
$ws.Range("K73").Value = 13999.9995
$ws.Range("I73").Value = 4666.6665
$ws.Range("M73").Value = -12907.9995
$ws.Range("N73").Value = -19168.9995
$ws.Range("L73").Value = 16984.9995
$ws.Range("H73").Value = 5330
$ws.Range("J73").Value = 5661.6665
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("L46").Value = 39496.668
$ws.Range("H46").Value = 28275.555
$ws.Range("J46").Value = 39496.668
$ws.Range("N46").Value = -39808.668
$ws.Range("H97").Value = 1814.1
$ws.Range("K97").Value = 1814.1
$ws.Range("M97").Value = -1318.1
$ws.Range("I97").Value = 1814.1
$ws.Range("N106").Value = -47524
$ws.Range("L106").Value = 45000
$ws.Range("H106").Value = 45000
$ws.Range("J106").Value = 45000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3259.7693
$ws.Range("I113").Value = 3259.7693
$ws.Range("N113").ClearContents()
$ws.Range("L113").Value = 0
$ws.Range("H113").Value = 3259.7693
$ws.Range("M113").Value = -1089.7693
$ws.Range("I132").Value = 3493.75
$ws.Range("H132").Value = 4558.5
$ws.Range("M132").Value = -7951.25
$ws.Range("K132").Value = 10481.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4590
$ws.Range("M22").Value = -3370.6667
$ws.Range("K22").Value = 3665.6667
$ws.Range("J22").Value = 4000
$ws.Range("H22").Value = 3799.4
$ws.Range("I22").Value = 3665.6667
$ws.Range("I27").Value = 3665.6667
$ws.Range("L27").Value = 4000
$ws.Range("J27").Value = 4000
$ws.Range("N27").Value = -4214
$ws.Range("M27").Value = -3558.6667
$ws.Range("H27").Value = 3799.4
$ws.Range("K27").Value = 3665.6667
$ws.Range("K93").Value = 974.3889
$ws.Range("N93").Value = -4196
$ws.Range("J93").Value = 1700
$ws.Range("H93").Value = 1012.5789
$ws.Range("I93").Value = 974.3889
$ws.Range("L93").Value = 1700
$ws.Range("M93").Value = 273.6111
$ws.Range("J109").Value = 0
$ws.Range("H109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("L109").Value = 0
$ws.Range("H122").Value = 391567.34
$ws.Range("M122").Value = -1445529.4
$ws.Range("I122").Value = 482659.8
$ws.Range("K122").Value = 1447979.4
$ws.Range("I132").Value = 4119.154
$ws.Range("H132").Value = 4139.2856
$ws.Range("M132").Value = -9827.462000000001
$ws.Range("K132").Value = 12357.462
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2232.5652
$ws.Range("M122").Value = -2962.6
$ws.Range("I122").Value = 1804.2
$ws.Range("N122").Value = -14007.25
$ws.Range("K122").Value = 5412.6
$ws.Range("J122").Value = 3035.75
$ws.Range("L122").Value = 9107.25
$ws.Range("J125").Value = 29400
$ws.Range("N125").Value = -39240
$ws.Range("H125").Value = 29400
$ws.Range("L125").Value = 29400
$ws.Range("J135").Value = 83916.664
$ws.Range("N135").Value = -94056.664
$ws.Range("L135").Value = 83916.664
